$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 66.388885
$ws.Range("I33").Value = 66.4375
$ws.Range("K33").Value = 66.4375
$ws.Range("M33").Value = 162.5625
$ws.Range("H40").Value = 1736.1538
$ws.Range("I40").Value = 1131.6666
$ws.Range("J40").Value = 2254.2856
$ws.Range("K40").Value = 1131.6666
$ws.Range("L40").Value = 2254.2856
$ws.Range("M40").Value = -956.6666
$ws.Range("N40").Value = -2604.2856
$ws.Range("H51").Value = 3083.5
$ws.Range("I51").Value = 3749.5
$ws.Range("J51").Value = 2750.5
$ws.Range("K51").Value = 3749.5
$ws.Range("L51").Value = 2750.5
$ws.Range("M51").Value = -3265.5
$ws.Range("N51").Value = -3718.5
$ws.Range("H55").Value = 333.10526
$ws.Range("I55").Value = 873.3333
$ws.Range("J55").Value = 83.76922999999999
$ws.Range("K55").Value = 873.3333
$ws.Range("L55").Value = 83.76922999999999
$ws.Range("M55").Value = -659.3333
$ws.Range("N55").Value = -511.76923
$ws.Range("H88").Value = 1211.7222
$ws.Range("I88").Value = 923.2
$ws.Range("J88").Value = 1322.6923
$ws.Range("K88").Value = 923.2
$ws.Range("L88").Value = 1322.6923
$ws.Range("M88").Value = -517.2
$ws.Range("N88").Value = -2134.6923
$ws.Range("H91").Value = 1211.7222
$ws.Range("I91").Value = 923.2
$ws.Range("J91").Value = 1322.6923
$ws.Range("K91").Value = 923.2
$ws.Range("L91").Value = 1322.6923
$ws.Range("M91").Value = 480.8
$ws.Range("N91").Value = -4130.6923
$ws.Range("H98").Value = 354.77777
$ws.Range("I98").Value = 369.7647
$ws.Range("K98").Value = 369.7647
$ws.Range("M98").Value = 1128.2353
$ws.Range("H111").Value = 4432.25
$ws.Range("I111").Value = 4432.25
$ws.Range("K111").Value = 13296.75
$ws.Range("M111").Value = -10229.75
$ws.Range("H122").Value = 354.77777
$ws.Range("I122").Value = 369.7647
$ws.Range("K122").Value = 1109.2941
$ws.Range("M122").Value = 1340.7059
$ws.Range("H129").Value = 186007.78
$ws.Range("J129").Value = 196927.84
$ws.Range("L129").Value = 590783.52
$ws.Range("N129").Value = -600783.52
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 87158.336
$ws.Range("I137").Value = 5666.5
$ws.Range("J137").Value = 168650.17
$ws.Range("K137").Value = 16999.5
$ws.Range("L137").Value = 505950.51
$ws.Range("M137").Value = -14449.5
$ws.Range("N137").Value = -511050.51
$ws.Range("H141").Value = 2466.25
$ws.Range("I141").Value = 1923.125
$ws.Range("J141").Value = 3552.5
$ws.Range("K141").Value = 5769.375
$ws.Range("L141").Value = 10657.5
$ws.Range("M141").Value = -589.375
$ws.Range("N141").Value = -21017.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 50002670
$ws.Range("I74").Value = 52634330
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 52634330
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -52633456
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 50002670
$ws.Range("I77").Value = 52634330
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 263171650
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -263167282
$ws.Range("N77").Value = -14736
$ws.Range("H122").Value = 2123.8635
$ws.Range("I122").Value = 2228.5
$ws.Range("K122").Value = 6685.5
$ws.Range("M122").Value = -4235.5
$ws.Range("H124").Value = 14534.75
$ws.Range("J124").Value = 14534.75
$ws.Range("L124").Value = 14534.75
$ws.Range("N124").Value = -24354.75
$ws.Range("H125").Value = 33943
$ws.Range("J125").Value = 33943
$ws.Range("L125").Value = 33943
$ws.Range("N125").Value = -43783
$ws.Range("H132").Value = 12345.766
$ws.Range("I132").Value = 1526.4286
$ws.Range("K132").Value = 4579.2858
$ws.Range("M132").Value = -2049.2858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H86").Value = 1726.1282
$ws.Range("J86").Value = 2551.4
$ws.Range("L86").Value = 2551.4
$ws.Range("N86").Value = -4797.4
$ws.Range("H89").Value = 1726.1282
$ws.Range("J89").Value = 2551.4
$ws.Range("L89").Value = 12757
$ws.Range("N89").Value = -23989
$ws.Range("H107").Value = 1699.5555
$ws.Range("I107").Value = 826.1429000000001
$ws.Range("K107").Value = 826.1429000000001
$ws.Range("M107").Value = 1093.8571
$ws.Range("H134").Value = 24637.588
$ws.Range("I134").Value = 30284.459
$ws.Range("K134").Value = 90853.37699999999
$ws.Range("M134").Value = -88318.37699999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1416
$ws.Range("I16").Value = 1520
$ws.Range("K16").Value = 1520
$ws.Range("M16").Value = -1233
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150
$ws.Range("H31").Value = 17628.947
$ws.Range("I31").Value = 22710.715
$ws.Range("J31").Value = 3400
$ws.Range("K31").Value = 22710.715
$ws.Range("L31").Value = 3400
$ws.Range("M31").Value = -22415.715
$ws.Range("N31").Value = -3990
$ws.Range("H34").Value = 17628.947
$ws.Range("I34").Value = 22710.715
$ws.Range("J34").Value = 3400
$ws.Range("K34").Value = 22710.715
$ws.Range("L34").Value = 3400
$ws.Range("M34").Value = -22508.715
$ws.Range("N34").Value = -3804
$ws.Range("H105").Value = 8929842
$ws.Range("I105").Value = 13889646
$ws.Range("J105").Value = 2194.4
$ws.Range("K105").Value = 13889646
$ws.Range("L105").Value = 2194.4
$ws.Range("M105").Value = -13887899
$ws.Range("N105").Value = -5688.4
$ws.Range("H113").Value = 1416
$ws.Range("I113").Value = 1520
$ws.Range("K113").Value = 1520
$ws.Range("M113").Value = 650
$ws.Range("H122").Value = 1855.6154
$ws.Range("I122").Value = 2517.5715
$ws.Range("K122").Value = 7552.7145
$ws.Range("M122").Value = -5102.7145
$ws.Range("H132").Value = 14942.975
$ws.Range("I132").Value = 19821.592
$ws.Range("J132").Value = 3966.0833
$ws.Range("K132").Value = 59464.776
$ws.Range("L132").Value = 11898.2499
$ws.Range("M132").Value = -56934.776
$ws.Range("N132").Value = -16958.2499

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 4447.4287
$ws.Range("I103").Value = 262
$ws.Range("J103").Value = 10028
$ws.Range("K103").Value = 786
$ws.Range("L103").Value = 30084
$ws.Range("M103").Value = 93
$ws.Range("N103").Value = -31842
$ws.Range("H131").Value = 748.61
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 748.61
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2245.83
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12325.83

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 22729444
$ws.Range("I102").Value = 27780228
$ws.Range("J102").Value = 908.5
$ws.Range("K102").Value = 27780228
$ws.Range("L102").Value = 908.5
$ws.Range("M102").Value = -27778606
$ws.Range("N102").Value = -4152.5
$ws.Range("H126").Value = 3178.0208
$ws.Range("I126").Value = 2332.7144
$ws.Range("J126").Value = 5453.846
$ws.Range("K126").Value = 6998.1432
$ws.Range("L126").Value = 16361.538
$ws.Range("M126").Value = -4528.1432
$ws.Range("N126").Value = -21301.538

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3162.8438
$ws.Range("J40").Value = 3892.45
$ws.Range("L40").Value = 3892.45
$ws.Range("N40").Value = -4164.45
$ws.Range("H46").Value = 1218.9
$ws.Range("I46").Value = 897.8
$ws.Range("J46").Value = 1540
$ws.Range("K46").Value = 897.8
$ws.Range("L46").Value = 1540
$ws.Range("M46").Value = -709.8
$ws.Range("N46").Value = -1916

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1156.12
$ws.Range("I132").Value = 645
$ws.Range("K132").Value = 1935
$ws.Range("M132").Value = 595
$ws.Range("H136").Value = 23257550
$ws.Range("I136").Value = 37038650
$ws.Range("K136").Value = 111115950
$ws.Range("M136").Value = -111113400
